$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77: extend task description, update hours, and let row height grow to fit the longer text
$ws.Range("D77").Value = "Got flagging system to work on topics result page.  Tested method to get visible unsuitable stories and started sketching out what to do with the admin."
$ws.Range("B77").Value = 5.5
$ws.Rows.Item(77).RowHeight = 30

# Row 82: swap the stale TODO note for the "NB" note that used to live on row 83
$ws.Range("D82").Value = "NB - the topic search results really aren't very useful without a profile link - I should do it if I can manage it!"

# Row 83 now duplicates what row 82 holds (it was moved up) - delete it so everything below shifts up one row
$ws.Rows.Item(83).Delete()

# Match the recorded selection state after the edit
$ws.Range("A82:XFD82").Select()
